$wb = $excel.ActiveWorkbook

# addCategory sheet: fix the "type" of the first category (was Digital, should be Physical)
$wsCategory = $wb.Worksheets.Item("addCategory")
$wsCategory.Range("D2").Value = "Physical"
$wsCategory.Range("A2").Select() | Out-Null

# addProduct sheet: update category name for the first product, move selection
$wsProduct = $wb.Worksheets.Item("addProduct")
$wsProduct.Range("B2").Value = "Hai Category 1"
$wsProduct.Range("C9").Select() | Out-Null

# Re-activate addProduct sheet (it was the tab-selected sheet)
$wsProduct.Activate() | Out-Null
